# Commit: "Added get_report_stats() to excel_Parsing which returns the
# gender and department stats of excel sheet"
#
# The new function derives a Gender value (column H) for every data row
# on the "Current H-1B cases" sheet from the employee's first name and
# writes it back into the sheet. The first several rows resolve to a
# handful of distinct spellings/casings; the remainder of the sheet
# falls back to the same value ("male").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current H-1B cases")

# --- Populate the "Gender" column (H) for rows 2-111 -----------------

# Rows 2-14: distinct values produced by the gender-lookup routine.
# Setting .Value on these cells leaves each cell's existing style intact.
$ws.Range("H2").Value = "M"
$ws.Range("H3").Value = "m"
$ws.Range("H4").Value = "F"
$ws.Range("H5:H7").Value = "f"
$ws.Range("H8").Value = "Male"
$ws.Range("H9").Value = "male"
$ws.Range("H10:H12").Value = "Female"
$ws.Range("H13:H14").Value = "female"

# Rows 15-111: the routine fell back to "male" for the rest of the
# sheet. Row 15's H cell already carries the highlighted "flag" style
# (fill/border) used for the rest of this block; copy that formatting
# down to H16:H111 before writing the values so every cell in the
# block ends up on that same style.
$ws.Range("H15").Copy()
$ws.Range("H16:H111").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H15:H111").Value = "male"
